# TC0001~TC0015 테스트 결과 입력
# - 이력(History) 시트: 새 이력 한 줄 추가 (row 11)
# - 테스트케이스 시트: TC0001~TC0015 에 대한 Pass/Fail 결과, 테스트일, 이슈번호 입력

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # 이력
$ws2 = $wb.Worksheets.Item(2)   # 테스트케이스

# ---------------------------------------------------------------------------
# 1) 이력 시트 - 새 이력 행 추가 (2017.06.05 / 3 / TC-0001~TC-0015 테스트 / 정화인)
# ---------------------------------------------------------------------------
# Copy the date value/type from the row above (A10) so "2017.06.05" stays
# plain text instead of being auto-parsed into a date serial number.
$ws1.Range("A10").Copy()
$ws1.Range("A11").PasteSpecial(-4163)   # xlPasteValues
$excel.CutCopyMode = $false

$ws1.Range("B11").Value = 3
$ws1.Range("C11").Value = "TC-0001~TC-0015 테스트"
$ws1.Range("D11").Value = "정화인"

# ---------------------------------------------------------------------------
# 2) 테스트케이스 시트 - TC0001~TC0015 실행 결과 입력 (행 9 ~ 23)
# ---------------------------------------------------------------------------
# Reuse the already-existing date style (numFmtId 14, same border/alignment
# as column J uses further down the sheet) by copying formats from J24 into
# J9:J23, then fill in the date value (2017-06-05 -> serial 42891).
$ws2.Range("J24").Copy()
$ws2.Range("J9:J23").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

for ($r = 9; $r -le 23; $r++) {
    $ws2.Cells.Item($r, 10).Value = 42891   # J column: 2017-06-05
}

# Rows that passed testing
$passRows = 9,10,11,12,13,15,17,19,20,21,22,23
foreach ($r in $passRows) {
    $ws2.Cells.Item($r, 9).Value = "Pass"
}

# Row 14 - Fail: 등록되지 않은 과목 수정 불가에 대한 에러 메시지가 출력되지않음.
$ws2.Range("I14").Value = "Fail " + [char]10 + "등록되지 않은 과목 수정 불가에 대한 에러 메시지가 출력되지않음."
$ws2.Range("I14").WrapText = $true
$ws2.Range("K14").Value = "#939"

# Row 16 - Fail: 빈 공간 선택에 대한 에러 메시지가 출력되지 않음
$ws2.Range("I16").Value = "Fail" + [char]10 + "빈 공간 선택에 대한 에러 메시지가 출력되지 않음"
$ws2.Range("I16").WrapText = $true
$ws2.Range("K16").Value = "#940"

# Row 18 - Fail: 항목명 미입력에 대한 에러 메시지가 출력되지않음
$ws2.Range("I18").Value = "Fail" + [char]10 + "항목명 미입력에 대한 에러 메시지가 출력되지않음"
$ws2.Range("I18").WrapText = $true
$ws2.Range("K18").Value = "#941"

# ---------------------------------------------------------------------------
# 3) 마지막 작업 화면 - 테스트케이스 시트의 마지막 입력 셀(K18)이 선택된 채로 저장
# ---------------------------------------------------------------------------
[void]$ws1.Range("D9").Select()
[void]$ws2.Activate()
[void]$ws2.Range("K18").Select()
